$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.323.36"
$ws.Range("E2").Value = "  -7.40%  "
$ws.Range("D3").Value = "3.276.65"
$ws.Range("E3").Value = "  -8.83%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "182.58"
$ws.Range("E5").Value = "  -11.75%  "
$ws.Range("D6").Value = "516.07"
$ws.Range("E6").Value = "  -9.07%  "
$ws.Range("D7").Value = "0.594"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("D8").Value = "3.270.43"
$ws.Range("E8").Value = "  -8.87%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "0.623"
$ws.Range("E10").Value = "  -8.35%  "
$ws.Range("D11").Value = "59.63"
$ws.Range("E11").Value = "  -6.78%  "
$ws.Range("D12").Value = "0.132"
$ws.Range("E12").Value = "  -10.94%  "
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  -8.42%  "
$ws.Range("D14").Value = "9.16"
$ws.Range("E14").Value = "  -9.17%  "
$ws.Range("D15").Value = "3.791.45"
$ws.Range("E15").Value = "  -8.80%  "
$ws.Range("E16").Value = "  -4.99%  "
$ws.Range("D17").Value = "3.272.10"
$ws.Range("E17").Value = "  -9.11%  "
$ws.Range("D18").Value = "17.47"
$ws.Range("E18").Value = "  -8.87%  "
$ws.Range("D19").Value = "63.152.17"
$ws.Range("E19").Value = "  -7.34%  "
$ws.Range("D20").Value = "11.00"
$ws.Range("E20").Value = "  -9.96%  "
$ws.Range("D21").Value = "0.949"
$ws.Range("E21").Value = "  -10.79%  "
$ws.Range("D22").Value = "372.55"
$ws.Range("E22").Value = "  -7.83%  "
$ws.Range("D23").Value = "11.38"
$ws.Range("E23").Value = "  -7.65%  "
$ws.Range("D24").Value = "80.29"
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("D25").Value = "3.66"
$ws.Range("E25").Value = "  -11.69%  "
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").Value = "5.92"
$ws.Range("E27").Value = "  -3.44%  "
$ws.Range("D28").Value = "2.66"
$ws.Range("E28").Value = "  -8.03%  "
$ws.Range("D29").Value = "11.39"
$ws.Range("E29").Value = "  -8.66%  "
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  -8.07%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "28.46"
$ws.Range("E31").Value = "  -9.83%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "644.33"
$ws.Range("E32").Value = "  -14.36%  "
$ws.Range("D33").Value = "6.74"
$ws.Range("E33").Value = "  -11.66%  "
$ws.Range("D34").Value = "11.22"
$ws.Range("E34").Value = "  -7.39%  "
$ws.Range("D35").Value = "59.21"
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("D36").Value = "0.105"
$ws.Range("E36").Value = "  -7.35%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "0.391"
$ws.Range("E38").Value = "  -8.24%  "
$ws.Range("D39").Value = "36.25"
$ws.Range("E39").Value = "  -12.74%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").Value = "  -5.14%  "
$ws.Range("D42").Value = "2.944.66"
$ws.Range("E42").Value = "  -7.29%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0663"
$ws.Range("E43").Value = "  -10.60%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  -6.11%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").Value = "  -17.40%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  +6.86%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.60"
$ws.Range("E47").Value = "  -6.31%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0389"
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.124"
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "2.46"
$ws.Range("E51").Value = "  -23.66%  "
